$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header area updates
$ws.Range("E11").Value = 343200   # Valor Mora (was 448539)
$ws.Range("C13").Value = 1        # Cant. Trabajadores (was 2)
$ws.Range("F13").Value = 11       # Cant. Periodos (was 14)

# 2) Re-order the worker detail rows (16-26) so the "Periodo Mora" column
#    goes from descending (1807..1709) to ascending (1709..1807), keeping
#    the rest of each row's data (Tipo Doc, N Doc, Nombre, Valor Mora,
#    Salario Basico) untouched since they are identical across the block.
$periods = @("1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# 3) Row 26 becomes the new last row of the block (after the next block is
#    removed below), so give it the "closing" bottom-border formatting that
#    used to belong to row 29 (the last row of the second, now-removed,
#    worker block).
$ws.Range("B29:J29").Copy() | Out-Null
$ws.Range("B26:J26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Remove the second worker block (rows 27-29, PEDRO ANTONIO URRIOLA LAMBIS)
$ws.Rows("27:29").Delete() | Out-Null

# 5) Remove the now-empty gap rows between the table and the signature
#    block: originally rows 30-33 were blank and the signature rows were
#    34-35; after removing the 3 rows above they are 31-32 and need no
#    further action -- nothing else occupies rows 27-30 any more.
